$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.413.50"
$ws.Range("E2").Value = "  -1.44%  "

$ws.Range("D3").Value = "2.282.32"
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "303.41"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").Value = "94.90"
$ws.Range("E6").Value = "  -3.37%  "

$ws.Range("D7").Value = "0.503"
$ws.Range("E7").Value = "  -3.07%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -3.52%  "

$ws.Range("D10").Value = "34.89"
$ws.Range("E10").Value = "  -3.99%  "

$ws.Range("E11").Value = "  -1.53%  "

$ws.Range("E12").Value = "  +1.54%  "

$ws.Range("D13").Value = "17.95"
$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("D14").Value = "6.78"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").Value = "2.633.85"
$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").Value = "2.294.79"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("E17").Value = "  -1.81%  "

$ws.Range("D18").Value = "42.327.02"
$ws.Range("E18").Value = "  -1.56%  "

$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  -0.41%  "

$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  -2.61%  "

$ws.Range("D21").Value = "5.94"
$ws.Range("E21").Value = "  -3.08%  "

$ws.Range("D22").Value = "67.01"
$ws.Range("E22").Value = "  -2.04%  "

$ws.Range("D23").Value = "235.69"
$ws.Range("E23").Value = "  -2.71%  "

$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  +0.65%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").Value = "2.38"
$ws.Range("E26").Value = "  -2.26%  "

$ws.Range("D27").Value = "24.59"
$ws.Range("E27").Value = "  -2.43%  "

$ws.Range("E28").Value = "  +17.07%  "

$ws.Range("E29").Value = "  +0.47%  "

$ws.Range("D31").Value = "32.20"
$ws.Range("E31").Value = "  -3.35%  "

$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").Value = "17.59"
$ws.Range("E33").Value = "  -0.88%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "4.91"
$ws.Range("E34").Value = "  -2.32%  "

$ws.Range("D35").Value = "4.43"
$ws.Range("E35").Value = "  -6.98%  "

$ws.Range("E36").Value = "  -2.41%  "

$ws.Range("D37").Value = "0.0682"
$ws.Range("E37").Value = "  -1.07%  "

$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("E39").Value = "  -2.88%  "

$ws.Range("E40").Value = "  -2.52%  "

$ws.Range("E41").Value = "  -4.71%  "

$ws.Range("D42").Value = "1.982.89"
$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("E43").Value = "  -4.20%  "

$ws.Range("D44").Value = "10.01"
$ws.Range("E44").Value = "  -2.26%  "

$ws.Range("D45").Value = "17.81"
$ws.Range("E45").Value = "  +2.10%  "

$ws.Range("E46").Value = "  -7.69%  "

$ws.Range("D48").Value = "2.91"
$ws.Range("E48").Value = "  +4.94%  "

$ws.Range("D49").Value = "53.22"
$ws.Range("E49").Value = "  -0.84%  "

$ws.Range("D50").Value = "2.502.24"
$ws.Range("E50").Value = "  -1.11%  "

$ws.Range("D51").Value = "70.25"
$ws.Range("E51").Value = "  -3.67%  "

